$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: convert inline string values to real numbers
$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = 1000

# New rows 3-10 with numeric data
$data = @(
    @(10050, 9950),
    @(1000, 1000),
    @(1000, 1000),
    @(2000, 1000),
    @(2000, 1000),
    @(2000, 1000),
    @(1000, 1000),
    @(5000, 5000)
)

$row = 3
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
